# Adds the next day's row (2025-03-23) to every price sheet in the workbook,
# mirroring the previous day's (row 21) value for each sheet, matching the
# "Updated Argent prices in Excel" commit.

$wb = $excel.ActiveWorkbook

$newDate = "2025-03-23"

# Sheet name -> new price value for column B, row 22
$updates = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.19"
    "Cell Topcon 183mm"          = "0.298"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,399"
    "Silver Busbar front-side"   = "8,083"
    "Silver finger front-side"   = "8,133"
    "USD_CNY"                    = "7.2717"
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 22
    $dateCell = $ws.Cells.Item($row, 1)
    $priceCell = $ws.Cells.Item($row, 2)
    $rowRange = $ws.Range($dateCell, $priceCell)

    # Force the cells to be treated as plain text so values such as "40" or
    # the date string aren't auto-converted into numbers/dates by Excel.
    $rowRange.NumberFormat = "@"

    $dateCell.Value = $newDate
    $priceCell.Value = $updates[$sheetName]

    # Remove the temporary text formatting so the new cells keep the default
    # (unstyled) look used by all the other rows in the sheet.
    $rowRange.ClearFormats()
}
